$d = $word.ActiveDocument

$d.Content.Find.Execute("492×9=", $true, $false, $false, $false, $false, $true, 1, $false, "985×8=", 2)
$d.Content.Find.Execute("246×8=", $true, $false, $false, $false, $false, $true, 1, $false, "849×8=", 2)
$d.Content.Find.Execute("829×8=", $true, $false, $false, $false, $false, $true, 1, $false, "581×6=", 2)
$d.Content.Find.Execute("282×6=", $true, $false, $false, $false, $false, $true, 1, $false, "309×7=", 2)
$d.Content.Find.Execute("342×5=", $true, $false, $false, $false, $false, $true, 1, $false, "830×7=", 2)
$d.Content.Find.Execute("909×4=", $true, $false, $false, $false, $false, $true, 1, $false, "224×2=", 2)
$d.Content.Find.Execute("648×8=", $true, $false, $false, $false, $false, $true, 1, $false, "574×3=", 2)
$d.Content.Find.Execute("285×2=", $true, $false, $false, $false, $false, $true, 1, $false, "863×2=", 2)
$d.Content.Find.Execute("623×9=", $true, $false, $false, $false, $false, $true, 1, $false, "371×9=", 2)
$d.Content.Find.Execute("349×4=", $true, $false, $false, $false, $false, $true, 1, $false, "183×2=", 2)
$d.Content.Find.Execute("233×2=", $true, $false, $false, $false, $false, $true, 1, $false, "139×3=", 2)
$d.Content.Find.Execute("622×2=", $true, $false, $false, $false, $false, $true, 1, $false, "554×3=", 2)
$d.Content.Find.Execute("660×6=", $true, $false, $false, $false, $false, $true, 1, $false, "985×2=", 2)
$d.Content.Find.Execute("374×5=", $true, $false, $false, $false, $false, $true, 1, $false, "520×3=", 2)
$d.Content.Find.Execute("364×6=", $true, $false, $false, $false, $false, $true, 1, $false, "778×7=", 2)
$d.Content.Find.Execute("529×8=", $true, $false, $false, $false, $false, $true, 1, $false, "824×3=", 2)
$d.Content.Find.Execute("420×3=", $true, $false, $false, $false, $false, $true, 1, $false, "981×3=", 2)
$d.Content.Find.Execute("932×9=", $true, $false, $false, $false, $false, $true, 1, $false, "951×2=", 2)
$d.Content.Find.Execute("694×9=", $true, $false, $false, $false, $false, $true, 1, $false, "612×7=", 2)
$d.Content.Find.Execute("736×7=", $true, $false, $false, $false, $false, $true, 1, $false, "200×7=", 2)
$d.Content.Find.Execute("465×3=", $true, $false, $false, $false, $false, $true, 1, $false, "622×7=", 2)
$d.Content.Find.Execute("274×2=", $true, $false, $false, $false, $false, $true, 1, $false, "375×6=", 2)
$d.Content.Find.Execute("122×8=", $true, $false, $false, $false, $false, $true, 1, $false, "436×4=", 2)
$d.Content.Find.Execute("197×6=", $true, $false, $false, $false, $false, $true, 1, $false, "188×9=", 2)
$d.Content.Find.Execute("529×9=", $true, $false, $false, $false, $false, $true, 1, $false, "631×5=", 2)
